# Performance Comparison.xlsx - Milestone1 update
# "added comments and file summary info.  Small optimizations to reduce processing time."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small optimizations: reduced processing time for the 10k/100k/1M rows ---
$ws.Range("P11").Value = 12
$ws.Range("P12").Value = 82
$ws.Range("P13").Value = 711

# --- New mini table (row 18) labelling the block below: "OpenCL" / "R" ---
$ws.Range("P18").Value = "OpenCL"
$ws.Range("V18").Value = "R"

# --- Row 19 : new run data, old N19 leftover removed ---
$ws.Range("N19").ClearContents()
$ws.Range("P19").Value = 12
$ws.Range("Q19").Value = 12
$ws.Range("R19").Value = 12
$ws.Range("V19").Value = 14
$ws.Range("W19").Value = 15
$ws.Range("X19").Value = 13

# --- Row 20 ---
$ws.Range("N20").ClearContents()
$ws.Range("O20").ClearContents()
$ws.Range("P20").Value = 13
$ws.Range("Q20").Value = 12
$ws.Range("R20").Value = 12
$ws.Range("V20").Value = 17
$ws.Range("W20").Value = 20
$ws.Range("X20").Value = 16

# --- Row 21 ---
$ws.Range("N21").ClearContents()
$ws.Range("O21").ClearContents()
$ws.Range("P21").Value = 12
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = 15
$ws.Range("V21").Value = 77
$ws.Range("W21").Value = 79
$ws.Range("X21").Value = 75

# --- Row 22 ---
$ws.Range("N22").ClearContents()
$ws.Range("O22").ClearContents()
$ws.Range("P22").Value = 82
$ws.Range("Q22").Value = 81
$ws.Range("R22").Value = 81
$ws.Range("V22").Value = 836
$ws.Range("W22").Value = 869
$ws.Range("X22").Value = 748

# --- Row 23 ---
$ws.Range("N23").ClearContents()
$ws.Range("O23").ClearContents()
$ws.Range("P23").Value = 711
$ws.Range("Q23").Value = 716
$ws.Range("R23").Value = 711
$ws.Range("V23").Value = 8074
$ws.Range("W23").Value = 8259
$ws.Range("X23").Value = 7755

# --- Row 24 : removed entirely (K25 block now starts right after row 23) ---
$ws.Range("N24").ClearContents()
$ws.Range("P24").ClearContents()
$ws.Range("Q24").ClearContents()

# --- Row 25 : only K25 (=10000/250) remains ---
$ws.Range("N25").ClearContents()
$ws.Range("P25").ClearContents()
$ws.Range("Q25").ClearContents()

# --- Rows 26-27 : removed entirely ---
$ws.Range("N26").ClearContents()
$ws.Range("P26").ClearContents()
$ws.Range("Q26").ClearContents()

$ws.Range("N27").ClearContents()
$ws.Range("P27").ClearContents()
$ws.Range("Q27").ClearContents()

# --- Update the view: scroll so column L is leftmost, select V19 ---
$excel.Goto($ws.Range("L1"), $true) | Out-Null
$ws.Range("V19").Select() | Out-Null
